$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.235.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.856.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7056'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08005'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3017'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08176'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.180'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7034'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.221.12'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007941'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.790'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.84'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.090.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.458'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.82%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.89'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  -0.95%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1429'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.916'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.428'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.20%  '
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.356'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.013'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05185'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("E35").Value = '  -1.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7127'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9986'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.648'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.724'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9426'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.132.25'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.936'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4250'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5310'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.759'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.985.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.166'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.27%  '
